$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.013453483581543
$ws.Range("B1").Value = 2.15196418762207
$ws.Range("C1").Value = 7.952270984649658
$ws.Range("D1").Value = 0.9781984090805054
$ws.Range("E1").Value = 0.5209988355636597
